# Thesis proposal glossary: add a new "Free word order language" /
# "زبان با ترتیب آزاد کلمات" row right after the existing
# "Sequence labeling" / "برچسب زنی دنباله (توالی)" row, at the end of
# the two-column English/Persian glossary table.

$d = $word.ActiveDocument

# Locate the "Sequence labeling" row by its text (more robust than a
# hard-coded row index) and resolve the table/row/cell it lives in.
$findRng = $d.Content
$findRng.Find.Execute("Sequence labeling") | Out-Null

$anchorCell = $findRng.Cells.Item(1)
$table = $anchorCell.Tables.Item(1)
$seqRow = $table.Rows.Item($anchorCell.RowIndex)

# Insert the new row immediately after the "Sequence labeling" row
# (i.e. before the row that currently follows it).
$insertBeforeRow = $table.Rows.Item($seqRow.Index + 1)
$newRow = $table.Rows.Add($insertBeforeRow)

# Left column: English term (left-to-right, default language run).
$enCell = $newRow.Cells.Item(1)
$enCell.Range.Text = "Free word order language "

# Right column: Persian translation (right-to-left, Farsi run).
$faCell = $newRow.Cells.Item(2)
$faCell.Range.Text = "زبان با ترتیب آزاد کلمات "

Write-Host "Row count now:" $table.Rows.Count
